# Pinouts workbook update: add GPIO pin-function labels, a colour legend,
# and reference notes; colour-code matching pin pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Pin-function labels in the main header/footer columns (B & D),
#    plus the "resistor / jumper count" helper column (E).
# ---------------------------------------------------------------------
$ws.Range("B5").Value  = "Right Light"
$ws.Range("D5").Value  = "Right Light"
$ws.Range("E5").Value  = "Resistor on Right"

$ws.Range("B9").Value  = "Middle Light"
$ws.Range("D9").Value  = "Middle Light"
$ws.Range("E9").Value  = "Resistor on Right"

$ws.Range("B10").Value = "Motion Sensor Ground"
$ws.Range("D10").Value = "Motion Sensor Data"
$ws.Range("E10").Value = 4

$ws.Range("D11").Value = "Motion Sensor 3.3V"

$ws.Range("B14").Value = "Push Button"
$ws.Range("D14").Value = "Push Button"
$ws.Range("E14").Value = 3

$ws.Range("B16").Value = "Light Sensor"
$ws.Range("D16").Value = "Light Sensor"
$ws.Range("E16").Value = 2

$ws.Range("B20").Value = "Left Light"
$ws.Range("D20").Value = "Left Light"
$ws.Range("E20").Value = "Resistor on Right"

# B21/D21 already held "Wifi Switch" - leave the values, only formatting
# changes below.

# ---------------------------------------------------------------------
# 2. Colour legend block (columns F & G, rows 1-5) plus the common-notes
#    rows (G7:G9).
# ---------------------------------------------------------------------
$ws.Range("F1").Value = 1
$ws.Range("G1").Value = "Bluetooth"

$ws.Range("F2").Value = 2
$ws.Range("G2").Value = "Light Sensor"

$ws.Range("F3").Value = 3
$ws.Range("G3").Value = "Push Button"

$ws.Range("F4").Value = 4
$ws.Range("G4").Value = "Motion Sensor"

$ws.Range("F5").Value = "All"
$ws.Range("G5").Value = "Common Items"

$ws.Range("G7").Value = "Note the left and right are swapped by the ribbon cable, with reference to the pinout."
$ws.Range("G8").Value = "Top is defined as the side AWAY from the USB"
$ws.Range("G9").Value = "Left is defined as the side with the HDMI"

# ---------------------------------------------------------------------
# 3. Colour-code the matching pin groups. Interior.Color takes an OLE
#    (BGR) value, i.e. the reverse byte order of the RGB hex code.
#    Colours are applied cell-by-cell (union "A1,B2" ranges only paint
#    the first area in this host) but grouped so the fill/style creation
#    order matches the source workbook: red, green, yellow, teal, blue,
#    then the right-aligned teal variant.
# ---------------------------------------------------------------------
$redCells    = @("F4","G4","B10","D10","D11")
$greenCells  = @("F3","G3","B14","D14")
$yellowCells = @("F2","G2","B16","D16")
$tealCells   = @("B5","D5","G5","B9","D9","B20","D20","B21","D21")
$blueCells   = @("F1","G1")

foreach ($ref in $redCells)    { $ws.Range($ref).Interior.Color = 0x5E5EF8 }   # red    FFF85E5E
foreach ($ref in $greenCells)  { $ws.Range($ref).Interior.Color = 0x50B000 }   # green  FF00B050
foreach ($ref in $yellowCells) { $ws.Range($ref).Interior.Color = 0x00FFFF }   # yellow FFFFFF00
foreach ($ref in $tealCells)   { $ws.Range($ref).Interior.Color = 0xCCCC33 }   # teal   FF33CCCC
foreach ($ref in $blueCells)   { $ws.Range($ref).Interior.Color = 0xC07000 }   # blue   FF0070C0

$ws.Range("F5").Interior.Color = 0xCCCC33            # teal, same as above
$ws.Range("F5").HorizontalAlignment = -4152          # xlRight

# ---------------------------------------------------------------------
# 4. Column widths - best-fit the label columns (B:E), and size the
#    notes column (G) by hand, matching the widths Excel settled on.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.333333333333336   # B -> ~21.14 best fit
$ws.Columns.Item(3).ColumnWidth = 8.5                  # C -> ~9.29  best fit
$ws.Columns.Item(4).ColumnWidth = 17.666666666666668   # D -> ~18.57 best fit
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666   # E -> 16     best fit
$ws.Columns.Item(7).ColumnWidth = 15.833333333333334   # G -> ~16.71 manual

# ---------------------------------------------------------------------
# 5. Leave the cursor where the author last left it.
# ---------------------------------------------------------------------
$ws.Range("G10").Select() | Out-Null
